# The edit reorders the observation records held in rows 2-11 of the
# "Artfynd" sheet: each row's data (Id, Antal, coordinates, timestamps,
# comments, ...) is replaced by the data that used to live in a different
# row, i.e. the ten records are permuted while staying on rows 2-11.
#
# Strategy: snapshot the current per-row values for every column that can
# differ between records, then write them back out in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based index) whose content can vary row-to-row for these
# records. Columns not listed here (C, L, N, P, T, U, V, W, Y, AA, AD, AE,
# AG, AT, AW, AX, AY, ...) are identical across rows 2-11 already, so they
# do not need to be touched.
$cols = @(1, 2, 4, 5, 6, 7, 8, 9, 10, 11, 17, 18, 19, 26, 28, 29)

$firstRow = 2
$lastRow = 11

# 1. Snapshot current values for rows 2..11.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# 2. New row order: row -> source row to copy from.
$mapping = @{
    2  = 11
    3  = 5
    4  = 2
    5  = 9
    6  = 4
    7  = 10
    8  = 6
    9  = 7
    10 = 3
    11 = 8
}

# 3. Write the permuted data back.
# Column I ("Antal") stores numeric-looking text (e.g. "20") as a genuine
# text value in the source file, not a number. Plain `.Value` assignment
# would let Excel auto-coerce a digit string back into a number, so that
# column is forced to text format first to keep it a string.
$antalCol = 9

foreach ($destRow in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $val = $srcData[$c]
        if ($null -eq $val) {
            $val = ""
        }
        $cell = $ws.Cells.Item($destRow, $c)
        if ($c -eq $antalCol -and $val -ne "") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}
